$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = 3517481.52
$ws.Range("C7").Value = -20.83232783672056
$ws.Range("D7").Value = 3073
$ws.Range("E7").Value = 3073
$ws.Range("F7").Value = 1144.640911161731
$ws.Range("G7").Value = 22.01044431021526
